$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows for the 5 part numbers that were removed (in descending
# order so earlier row indices remain valid while deleting).
$ws.Rows("78:78").Delete()
$ws.Rows("77:77").Delete()
$ws.Rows("46:46").Delete()
$ws.Rows("42:42").Delete()
$ws.Rows("23:23").Delete()

# The used range shrank from 78 to 73 data rows; refresh the AutoFilter to
# match, and keep the hidden _xlnm._FilterDatabase defined name in sync too.
$ws.AutoFilterMode = $false
$ws.Range("A1:L73").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=données!`$A`$1:`$L`$73"
    }
}

# Restore the view state recorded after the edit.
$ws.Range("B74:B75").Select()
$excel.ActiveWindow.ScrollColumn = 2
